# Apply updated odds values to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "N2"  = 1.02
    "P2"  = 1.02
    "R2"  = 1.02
    "N3"  = 1.02
    "P3"  = 1.02
    "R3"  = 1.02
    "P4"  = 1.46
    "L6"  = 1.46
    "M6"  = 1.08
    "T6"  = 1.87
    "S7"  = 1.78
    "R9"  = 1.32
    "P12" = 1.79
    "AN12" = 18
    "I13" = 3.3
    "L13" = 1.41
    "N15" = 5.3
    "S18" = 2.08
    "AO18" = 8
    "G19" = 1.64
    "L20" = 1.27
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
